$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update quantity and price
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 4199.7

# Row 3: the part number text now resolves to "7114168" after the
# removal of two shared-string entries; update quantity and price too
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "7114168"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 1431

# Row 4: only the total remains, in column C; clear the old A4/B4/C4
# contents and move the (updated) total label into C4
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "Total: 5630.7"

# Rows 5 and 6 are removed entirely
$ws.Range("A5:C5").ClearContents()
$ws.Range("A6:C6").ClearContents()
